$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 61; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "'2024-09-26"
}
